# Generate Report for Handoff
# - Status "In Translation" -> "Ready for handoff" (Overview B2/C2, zh-cn C2, de-de C2)
# - Overview's "Latest Handoff Date" timestamp refreshed
# - zh-cn / de-de "Latest Handoff Datetime" timestamps refreshed

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns for both locales + latest handoff date
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-37-13 04:37:21"

# zh-cn detail sheet: Status + Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-13 04:37:17"

# de-de detail sheet: Status + Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-13 04:37:21"
